$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF" + $row)
    if ($cell.Value2 -eq "5-7-2013-14") {
        $cell.Value = "2014-05-07"
    }
}
